$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Datas" column header
$ws.Range("C1").Value = "Datas"

# Date values (stored as date serials, matching 2022-01-22 18:00 and 2022-01-21 15:00)
$ws.Range("C2").Value = 44583.75
$ws.Range("C3").Value = 44582.625

# Apply a short-date number format to C2, then propagate the exact same style to C3
# via copy/paste-special so both cells share a single style entry.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# Resize column C to fit its new date content.
$ws.Columns("C").ColumnWidth = 9.85

# Move the active selection to the new header cell.
$null = $ws.Range("C1").Select()
